# Auto-generated from diff: updates market-data-driven numeric columns (H-N)
# across multiple sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2178.8333
$ws.Range("I19").Value = 2830
$ws.Range("J19").Value = 1992.7858
$ws.Range("K19").Value = 2830
$ws.Range("L19").Value = 1992.7858
$ws.Range("M19").Value = -2655
$ws.Range("N19").Value = -2342.7858
$ws.Range("H76").Value = 23818512
$ws.Range("I76").Value = 55569452
$ws.Range("J76").Value = 5306.6665
$ws.Range("K76").Value = 55569452
$ws.Range("L76").Value = 5306.6665
$ws.Range("M76").Value = -55569137
$ws.Range("N76").Value = -5936.6665
$ws.Range("H79").Value = 23818512
$ws.Range("I79").Value = 55569452
$ws.Range("J79").Value = 5306.6665
$ws.Range("K79").Value = 55569452
$ws.Range("L79").Value = 5306.6665
$ws.Range("M79").Value = -55568360
$ws.Range("N79").Value = -7490.6665
$ws.Range("H116").Value = 43347.46
$ws.Range("I116").Value = 78353.92999999999
$ws.Range("J116").Value = 2506.5833
$ws.Range("K116").Value = 78353.92999999999
$ws.Range("L116").Value = 2506.5833
$ws.Range("M116").Value = -74911.92999999999
$ws.Range("N116").Value = -9390.5833
$ws.Range("H119").Value = 2566.6667
$ws.Range("I119").Value = 2450
$ws.Range("J119").Value = 2625
$ws.Range("K119").Value = 7350
$ws.Range("L119").Value = 7875
$ws.Range("M119").Value = -2512
$ws.Range("N119").Value = -17551
$ws.Range("H125").Value = 816.8570999999999
$ws.Range("I125").Value = 715.0769
$ws.Range("J125").Value = 982.25
$ws.Range("K125").Value = 6435.6921
$ws.Range("L125").Value = 8840.25
$ws.Range("M125").Value = -3975.6921
$ws.Range("N125").Value = -13760.25
$ws.Range("H132").Value = 5277.4546
$ws.Range("I132").Value = 1320.3704
$ws.Range("J132").Value = 23084.334
$ws.Range("K132").Value = 3961.1112
$ws.Range("L132").Value = 69253.00199999999
$ws.Range("M132").Value = -1431.1112
$ws.Range("N132").Value = -74313.00199999999
$ws.Range("H135").Value = 202.8
$ws.Range("I135").Value = 202.8
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 1825.2
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = 709.8
$ws.Range("H137").Value = 2625626.2
$ws.Range("I137").Value = 3517455.8
$ws.Range("J137").Value = 1818733
$ws.Range("K137").Value = 10552367.4
$ws.Range("L137").Value = 5456199
$ws.Range("M137").Value = -10549817.4
$ws.Range("N137").Value = -5461299

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 653.7368
$ws.Range("I2").Value = 728.70966
$ws.Range("J2").Value = 321.7143
$ws.Range("K2").Value = 728.70966
$ws.Range("L2").Value = 321.7143
$ws.Range("M2").Value = -615.70966
$ws.Range("N2").Value = -547.7143
$ws.Range("H32").Value = 1638655.4
$ws.Range("I32").Value = 1638655.4
$ws.Range("K32").Value = 1638655.4
$ws.Range("M32").Value = -1638368.4
$ws.Range("H45").Value = 939.5
$ws.Range("I45").Value = 911.5
$ws.Range("K45").Value = 911.5
$ws.Range("M45").Value = -534.5
$ws.Range("H116").Value = 653.7368
$ws.Range("I116").Value = 728.70966
$ws.Range("J116").Value = 321.7143
$ws.Range("K116").Value = 728.70966
$ws.Range("L116").Value = 321.7143
$ws.Range("M116").Value = 1565.29034
$ws.Range("N116").Value = -4909.7143
$ws.Range("H132").Value = 36011.2
$ws.Range("I132").Value = 64340.625
$ws.Range("K132").Value = 193021.875
$ws.Range("M132").Value = -190491.875

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 653.7368
$ws.Range("I3").Value = 728.70966
$ws.Range("J3").Value = 321.7143
$ws.Range("K3").Value = 728.70966
$ws.Range("L3").Value = 321.7143
$ws.Range("M3").Value = -614.70966
$ws.Range("N3").Value = -549.7143
$ws.Range("H64").Value = 469.57574
$ws.Range("I64").Value = 472.94736
$ws.Range("J64").Value = 465
$ws.Range("K64").Value = 472.94736
$ws.Range("L64").Value = 465
$ws.Range("M64").Value = -247.94736
$ws.Range("N64").Value = -915
$ws.Range("H67").Value = 469.57574
$ws.Range("I67").Value = 472.94736
$ws.Range("J67").Value = 465
$ws.Range("K67").Value = 472.94736
$ws.Range("L67").Value = 465
$ws.Range("M67").Value = 307.05264
$ws.Range("N67").Value = -2025
$ws.Range("H86").Value = 8852.777
$ws.Range("I86").Value = 16075
$ws.Range("J86").Value = 3075
$ws.Range("K86").Value = 16075
$ws.Range("L86").Value = 3075
$ws.Range("M86").Value = -14952
$ws.Range("N86").Value = -5321
$ws.Range("H89").Value = 8852.777
$ws.Range("I89").Value = 16075
$ws.Range("J89").Value = 3075
$ws.Range("K89").Value = 80375
$ws.Range("L89").Value = 15375
$ws.Range("M89").Value = -74759
$ws.Range("N89").Value = -26607
$ws.Range("H107").Value = 1306.3334
$ws.Range("I107").Value = 1299.95
$ws.Range("J107").Value = 1338.25
$ws.Range("K107").Value = 1299.95
$ws.Range("L107").Value = 1338.25
$ws.Range("M107").Value = 620.05
$ws.Range("N107").Value = -5178.25
$ws.Range("H134").Value = 8592.375
$ws.Range("I134").Value = 13024.385
$ws.Range("J134").Value = 3354.5454
$ws.Range("K134").Value = 39073.155
$ws.Range("L134").Value = 10063.6362
$ws.Range("M134").Value = -36538.155
$ws.Range("N134").Value = -15133.6362

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3310.0908
$ws.Range("I58").Value = 3460.1428
$ws.Range("J58").Value = 3047.5
$ws.Range("K58").Value = 3460.1428
$ws.Range("L58").Value = 3047.5
$ws.Range("M58").Value = -3257.1428
$ws.Range("N58").Value = -3453.5
$ws.Range("H107").Value = 412.65714
$ws.Range("I107").Value = 260.37036
$ws.Range("J107").Value = 926.625
$ws.Range("K107").Value = 260.37036
$ws.Range("L107").Value = 926.625
$ws.Range("M107").Value = 1659.62964
$ws.Range("N107").Value = -4766.625
$ws.Range("H132").Value = 16131063
$ws.Range("I132").Value = 23810810
$ws.Range("J132").Value = 3594.2
$ws.Range("K132").Value = 71432430
$ws.Range("L132").Value = 10782.6
$ws.Range("M132").Value = -71429900
$ws.Range("N132").Value = -15842.6
$ws.Range("H134").Value = 55557610
$ws.Range("I134").Value = 100001000
$ws.Range("J134").Value = 3361
$ws.Range("K134").Value = 300003000
$ws.Range("L134").Value = 10083
$ws.Range("M134").Value = -300000465
$ws.Range("N134").Value = -15153
$ws.Range("H136").Value = 3310.0908
$ws.Range("I136").Value = 3460.1428
$ws.Range("J136").Value = 3047.5
$ws.Range("K136").Value = 10380.4284
$ws.Range("L136").Value = 9142.5
$ws.Range("M136").Value = -7830.428400000001
$ws.Range("N136").Value = -14242.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 777.36365
$ws.Range("I5").Value = 460
$ws.Range("K5").Value = 1380
$ws.Range("M5").Value = -1268
$ws.Range("H68").Value = 959.7843
$ws.Range("I68").Value = 585.86365
$ws.Range("J68").Value = 1243.4482
$ws.Range("K68").Value = 1757.59095
$ws.Range("L68").Value = 3730.3446
$ws.Range("M68").Value = -946.59095
$ws.Range("N68").Value = -5352.3446
$ws.Range("H70").Value = 75696.57000000001
$ws.Range("I70").Value = 201598.4
$ws.Range("J70").Value = 5751.1113
$ws.Range("K70").Value = 604795.2
$ws.Range("L70").Value = 17253.3339
$ws.Range("M70").Value = -604480.2
$ws.Range("N70").Value = -17883.3339
$ws.Range("H71").Value = 959.7843
$ws.Range("I71").Value = 585.86365
$ws.Range("J71").Value = 1243.4482
$ws.Range("K71").Value = 5272.77285
$ws.Range("L71").Value = 11191.0338
$ws.Range("M71").Value = -1216.77285
$ws.Range("N71").Value = -19303.0338
$ws.Range("H73").Value = 75696.57000000001
$ws.Range("I73").Value = 201598.4
$ws.Range("J73").Value = 5751.1113
$ws.Range("K73").Value = 604795.2
$ws.Range("L73").Value = 17253.3339
$ws.Range("M73").Value = -603703.2
$ws.Range("N73").Value = -19437.3339
$ws.Range("H113").Value = 604.86664
$ws.Range("I113").Value = 590.34784
$ws.Range("J113").Value = 652.5714
$ws.Range("K113").Value = 1771.04352
$ws.Range("L113").Value = 1957.7142
$ws.Range("M113").Value = 398.9564799999998
$ws.Range("N113").Value = -6297.7142
$ws.Range("H122").Value = 422.95834
$ws.Range("I122").Value = 266
$ws.Range("J122").Value = 535.0714
$ws.Range("K122").Value = 2394
$ws.Range("L122").Value = 4815.6426
$ws.Range("M122").Value = 56
$ws.Range("N122").Value = -9715.642599999999
$ws.Range("H135").Value = 777.36365
$ws.Range("I135").Value = 460
$ws.Range("K135").Value = 4140
$ws.Range("M135").Value = -1605

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 50000
$ws.Range("J74").Value = 50000
$ws.Range("L74").Value = 50000
$ws.Range("N74").Value = -51872
$ws.Range("H77").Value = 50000
$ws.Range("J77").Value = 50000
$ws.Range("L77").Value = 150000
$ws.Range("N77").Value = -159360
$ws.Range("H132").Value = 2910696.8
$ws.Range("I132").Value = 5438444.5
$ws.Range("J132").Value = 3786.8
$ws.Range("K132").Value = 16315333.5
$ws.Range("L132").Value = 11360.4
$ws.Range("M132").Value = -16312803.5
$ws.Range("N132").Value = -16420.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31872
$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99360

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 17534498
$ws.Range("I136").Value = 19078650
$ws.Range("J136").Value = 10094491
$ws.Range("K136").Value = 57235950
$ws.Range("L136").Value = 30283473
$ws.Range("M136").Value = -57233400
$ws.Range("N136").Value = -30288573

# --- Special case: ALC row 135 lost its N column cell (merged away) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N135").ClearContents()
